# Update the public EPEX Spot / Gaz / CO2 price workbook with the latest
# daily data point(s):
#   - "Prix Spot" sheet: new column CL ("11-sep") with 24 hourly prices
#   - "Gaz" sheet: new row 87 for 2025-09-09
#   - "CO2" sheet: new row 87 for 2025-09-09

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Prix Spot" sheet — append column CL (column 90) after CK (89)
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Clone the header cell's look (bold / centered / bordered) from the last
# existing header (CK1) onto the new header cell (CL1), then set its text.
$wsSpot.Range("CK1").Copy()
$wsSpot.Range("CL1").PasteSpecial(-4122)   # xlPasteFormats
$wsSpot.Range("CL1").Value = "11-sep"

# Hourly values for 11-sep, rows 2..25 (row 1 is the header).
$spotValues = @(14, 11.7, 10.53, 0, 0, 8.84, 9.640000000000001, 28.59, 37.08, 17.13, 0.03, -0.01, -0.01, -0.01, -0.01, -0.01, -0.01, 0, 7.05, 32.46, 25.85, 11.81, 15.97, 11.87)
for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 90).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------
# 2) "Gaz" sheet — append row 87 for 2025-09-09
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Write the date as literal text (matching the existing A-column cells,
# which are plain inline strings, not real Excel dates) by building it via
# a text formula in a scratch cell and pasting the computed value back —
# this avoids the automatic "looks like a date" conversion that a direct
# .Value assignment of "2025-09-09" would trigger.
$wsGaz.Range("ZZ1").Formula = "=""2025-09-09"""
$wsGaz.Range("ZZ1").Copy()
$wsGaz.Range("A87").PasteSpecial(-4163)    # xlPasteValues
$wsGaz.Range("ZZ1").ClearContents()

$wsGaz.Range("B87").Value = 32

# ---------------------------------------------------------------------
# 3) "CO2" sheet — append row 87 for 2025-09-09
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("ZZ1").Formula = "=""2025-09-09"""
$wsCo2.Range("ZZ1").Copy()
$wsCo2.Range("A87").PasteSpecial(-4163)    # xlPasteValues
$wsCo2.Range("ZZ1").ClearContents()

$wsCo2.Range("B87").Value = 75.8
